$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new 2024 Bledug Kuwu data row (row 62) ---
# Write the Comment (F) and Location (B) text cells first so the shared
# string table grows in the same order as the authoritative edit did
# ("Tingay, Pers. Comm." = index 74, "Bledug Kuwu, ..." = index 75).
$ws.Cells.Item(62, 6).Value = "Tingay, Pers. Comm."

$ws.Cells.Item(62, 2).Value = "Bledug Kuwu, Cangkring, East Java"
# Match the pasted-from-web text styling (dark grey, default Calibri 12)
# used for this new location cell.
$ws.Cells.Item(62, 2).Font.Color = 2236962

$ws.Cells.Item(62, 1).Value = 45373
$ws.Cells.Item(62, 1).NumberFormat = "m/d/yy"

$ws.Cells.Item(62, 3).Value = 6.4
$ws.Cells.Item(62, 4).Value = 195

# --- Widen column B so the longer location text fits ---
$ws.Columns.Item(2).ColumnWidth = 30

# --- Update the view state to match where the edit left the selection ---
[void]$ws.Range("B62").Select()
$excel.ActiveWindow.Zoom = 178
